$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '24.488.17'
$ws.Cells.Item(2, 5).Value = '  -1.36%  '
$ws.Cells.Item(3, 4).Value = '1.657.41'
$ws.Cells.Item(3, 5).Value = '  -2.92%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.004'
$ws.Cells.Item(4, 5).Value = '  +0.30%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '313.40'
$ws.Cells.Item(5, 5).Value = '  -0.51%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '1.005'
$ws.Cells.Item(6, 5).Value = '  +0.38%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.3928'
$ws.Cells.Item(7, 5).Value = '  -1.82%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3907'
$ws.Cells.Item(8, 5).Value = '  -3.29%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '1.012'
$ws.Cells.Item(9, 5).Value = '  +1.02%  '
$ws.Cells.Item(10, 2).Value = 'OKB'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '50.50'
$ws.Cells.Item(11, 2).Value = 'Polygon'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '1.381'
$ws.Cells.Item(11, 5).Value = '  -6.30%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.08562'
$ws.Cells.Item(12, 5).Value = '  -2.84%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '24.90'
$ws.Cells.Item(13, 5).Value = '  -5.45%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '7.221'
$ws.Cells.Item(14, 5).Value = '  -3.91%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.00001304'
$ws.Cells.Item(15, 5).Value = '  -2.74%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '7.615'
$ws.Cells.Item(16, 5).Value = '  -4.80%  '
$ws.Cells.Item(17, 4).Value = '1.664.56'
$ws.Cells.Item(17, 5).Value = '  -5.34%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '93.04'
$ws.Cells.Item(18, 5).Value = '  -2.61%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.07004'
$ws.Cells.Item(19, 5).Value = '  -2.30%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '21.26'
$ws.Cells.Item(20, 5).Value = '  +1.68%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.975'
$ws.Cells.Item(21, 5).Value = '  -4.27%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '1.005'
$ws.Cells.Item(22, 5).Value = '  +0.42%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '13.76'
$ws.Cells.Item(23, 5).Value = '  -4.73%  '
$ws.Cells.Item(24, 4).Value = '24.335.18'
$ws.Cells.Item(24, 5).Value = '  -1.91%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.349'
$ws.Cells.Item(25, 5).Value = '  -0.14%  '
$ws.Cells.Item(26, 5).Value = '  -4.40%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '22.76'
$ws.Cells.Item(27, 5).Value = '  -1.47%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '5.791'
$ws.Cells.Item(28, 5).Value = '  -9.39%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '158.24'
$ws.Cells.Item(29, 5).Value = '  -2.14%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '144.91'
$ws.Cells.Item(30, 5).Value = '  +0.79%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '8.327'
$ws.Cells.Item(31, 5).Value = '  +1.35%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '2.529'
$ws.Cells.Item(32, 5).Value = '  +11.04%  '
$ws.Cells.Item(33, 4).Value = '1.865.65'
$ws.Cells.Item(33, 5).Value = '  -0.62%  '
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.08201'
$ws.Cells.Item(34, 5).Value = '  -5.26%  '
$ws.Cells.Item(35, 2).Value = 'VeChain'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.03014'
$ws.Cells.Item(35, 5).Value = '  -5.83%  '
$ws.Cells.Item(36, 2).Value = 'ImmutableX'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.9946'
$ws.Cells.Item(36, 5).Value = '  -3.48%  '
$ws.Cells.Item(37, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '6.850'
$ws.Cells.Item(37, 5).Value = '  -6.11%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.2769'
$ws.Cells.Item(38, 5).Value = '  -3.03%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.09600'
$ws.Cells.Item(39, 5).Value = '  +1.67%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.499'
$ws.Cells.Item(40, 5).Value = '  +1.27%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '10.26'
$ws.Cells.Item(41, 5).Value = '  -4.45%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.7766'
$ws.Cells.Item(42, 5).Value = '  -7.64%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '13.29'
$ws.Cells.Item(43, 5).Value = '  -6.65%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '16.36'
$ws.Cells.Item(44, 5).Value = '  -6.13%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '2.528'
$ws.Cells.Item(45, 5).Value = '  -7.29%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.7000'
$ws.Cells.Item(46, 5).Value = '  -5.68%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '4.155'
$ws.Cells.Item(47, 5).Value = '  -1.55%  '
$ws.Cells.Item(48, 2).Value = 'Frax'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '1.006'
$ws.Cells.Item(48, 5).Value = '  +0.53%  '
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.08563'
$ws.Cells.Item(49, 5).Value = '  +2.14%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.306'
$ws.Cells.Item(50, 5).Value = '  -4.65%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '136.82'
$ws.Cells.Item(51, 5).Value = '  -2.76%  '
